$d = $word.ActiveDocument

# 1. Merge "Age and Gender" + " Analysis of customers coupon acceptance – " -- no
#    visible text change needed (Find/Replace can't alter run structure), but the
#    text content stays identical, so nothing to do here for plain text.

# 2. Change "Below age 21 show greater tendency to accept the coupons."
#    to "At age 21 it shows the greater tendency to accept the coupons."
$d.Content.Find.Execute("Below age 21 show greater tendency to accept the coupons.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "At age 21 it shows the greater tendency to accept the coupons.", 2)

# 3. Merge "Healthcare Support..." + "accepts" + " the least 54%." removing proofErr marks
#    Text content is identical already; proofErr removal isn't achievable via Find/Replace.
